# "update footer dan maps" - add two new candidate rows (LIA ETIKASARI,
# KURNIA AINUN) to the report sheet and widen columns A and C so the new
# longer values fit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 4: LIA ETIKASARI
$ws.Cells.Item(4, 1).Value = "LIA ETIKASARI"
$ws.Cells.Item(4, 2).Value = "Tidak Ingin Menyebutkan"
$ws.Cells.Item(4, 3).Value = "liaetikasari0826@gmail.com"
$ws.Cells.Item(4, 4).Value = 85158040206
$ws.Cells.Item(4, 5).Font.Name = "Calibri"
$ws.Cells.Item(4, 6).Font.Name = "Calibri"
$ws.Cells.Item(4, 7).Font.Name = "Calibri"
$ws.Cells.Item(4, 9).Value = "AKTIF"

# New row 5: KURNIA AINUN
$ws.Cells.Item(5, 1).Value = "KURNIA AINUN"
$ws.Cells.Item(5, 2).Value = "Tidak Ingin Menyebutkan"
$ws.Cells.Item(5, 3).Value = "etikasarilia26@gmail.com"
$ws.Cells.Item(5, 4).Value = 628885122711
$ws.Cells.Item(5, 5).Font.Name = "Calibri"
$ws.Cells.Item(5, 6).Font.Name = "Calibri"
$ws.Cells.Item(5, 7).Font.Name = "Calibri"
$ws.Cells.Item(5, 9).Value = "AKTIF"

# Widen column A (15 -> 16) and column C (24 -> 31) to fit the new data.
# Excel's ColumnWidth property is offset by 5/6 of a character from the
# raw OOXML <col width> value, so subtract that to land on exact widths.
$ws.Columns.Item(1).ColumnWidth = 16 - 5/6
$ws.Columns.Item(3).ColumnWidth = 31 - 5/6

Write-Host "applied"
